# Auto-generated cell update script
# Applies the numeric value changes to the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets
# as captured from the authoritative before/after workbook diff (scheduled
# profit-tracking data refresh).

$wb = $excel.ActiveWorkbook

$changes = @{
    "ALC" = @(
        @("H12", 0),
        @("I12", 0),
        @("J12", 0),
        @("K12", 0),
        @("L12", 0),
        @("M12", $null),
        @("N12", $null),
        @("H18", 1324.2858),
        @("H28", 2837.2307),
        @("I28", 2554.476),
        @("J28", 4024.8),
        @("K28", 2554.476),
        @("L28", 4024.8),
        @("M28", -2069.476),
        @("N28", -4994.8),
        @("H62", 18523352),
        @("I62", 22227022),
        @("K62", 22227022),
        @("M62", -22226398),
        @("H65", 18523352),
        @("I65", 22227022),
        @("K65", 111135110),
        @("M65", -111131990)
    )
    "ARM" = @(
        @("H32", 1994.73),
        @("I32", 1557.6818),
        @("J32", 5199.75),
        @("K32", 1557.6818),
        @("L32", 5199.75),
        @("M32", -1270.6818),
        @("N32", -5773.75)
    )
    "BSM" = @(
        @("H105", 83335020),
        @("I105", 90910744),
        @("K105", 90910744),
        @("M105", -90908997)
    )
    "CRP" = @(
        @("H22", 117046.836),
        @("J22", 175320.25),
        @("L22", 175320.25),
        @("N22", -176020.25),
        @("H31", 1636.95),
        @("I31", 1482.6666),
        @("J31", 2364.2856),
        @("K31", 1482.6666),
        @("L31", 2364.2856),
        @("M31", -1187.6666),
        @("N31", -2954.2856),
        @("H34", 1636.95),
        @("I34", 1482.6666),
        @("J34", 2364.2856),
        @("K34", 1482.6666),
        @("L34", 2364.2856),
        @("M34", -1280.6666),
        @("N34", -2768.2856),
        @("H58", 1567.5946),
        @("I58", 1221.8636),
        @("J58", 2074.6667),
        @("K58", 1221.8636),
        @("L58", 2074.6667),
        @("M58", -1018.8636),
        @("N58", -2480.6667),
        @("H122", 756.3889),
        @("I122", 713.375),
        @("J122", 1100.5),
        @("K122", 2140.125),
        @("L122", 3301.5),
        @("M122", 309.875),
        @("N122", -8201.5),
        @("H134", 35716330),
        @("I134", 2220.3333),
        @("J134", 100001720),
        @("K134", 6660.999899999999),
        @("L134", 300005160),
        @("M134", -4125.999899999999),
        @("N134", -300010230),
        @("H136", 1567.5946),
        @("I136", 1221.8636),
        @("J136", 2074.6667),
        @("K136", 3665.5908),
        @("L136", 6224.000100000001),
        @("M136", -1115.5908),
        @("N136", -11324.0001)
    )
    "CUL" = @(
        @("H14", 262.30768),
        @("I14", 262.30768),
        @("K14", 786.92304),
        @("M14", -613.92304),
        @("H40", 187.94444),
        @("I40", 86.63636),
        @("J40", 347.14285),
        @("K40", 346.54544),
        @("L40", 1388.5714),
        @("M40", -277.54544),
        @("N40", -1526.5714),
        @("H57", 0),
        @("I57", 0),
        @("K57", 0),
        @("M57", $null),
        @("H133", 4357.136),
        @("J133", 5193.1333),
        @("L133", 15579.3999),
        @("N133", -25699.3999),
        @("H134", 2679.0833),
        @("J134", 0),
        @("L134", 0),
        @("N134", $null),
        @("H137", 46884468),
        @("I137", 83335380),
        @("J137", 19009.428),
        @("K137", 250006140),
        @("L137", 57028.284),
        @("M137", -250001040),
        @("N137", -67228.284),
        @("H139", 2066.5),
        @("J139", 1582.3334),
        @("L139", 4747.0002),
        @("N139", -15027.0002),
        @("H140", 2679.1355),
        @("I140", 1919.5358),
        @("J140", 3365.2258),
        @("K140", 5758.607400000001),
        @("L140", 10095.6774),
        @("M140", -578.6074000000008),
        @("N140", -20455.6774),
        @("H141", 125003590),
        @("I141", 166668540),
        @("K141", 500005620),
        @("M141", -500000440)
    )
    "GSM" = @(
        @("H29", 501),
        @("I29", 501),
        @("J29", 0),
        @("K29", 501),
        @("L29", 0),
        @("M29", -211),
        @("N29", $null),
        @("H33", 0),
        @("J33", 0),
        @("L33", 0),
        @("N33", $null),
        @("H36", 0),
        @("I36", 0),
        @("K36", 0),
        @("M36", $null),
        @("H41", 500),
        @("I41", 0),
        @("J41", 500),
        @("K41", 0),
        @("L41", 500),
        @("M41", $null),
        @("N41", -1210),
        @("H46", 11916.667),
        @("I46", 2750),
        @("J46", 16500),
        @("K46", 2750),
        @("L46", 16500),
        @("M46", -2594),
        @("N46", -16812),
        @("H49", 23000),
        @("J49", 23000),
        @("L49", 23000),
        @("N49", -23368),
        @("H53", 4000),
        @("J53", 4000),
        @("L53", 4000),
        @("N53", -5262),
        @("H57", 23000),
        @("J57", 23000),
        @("L57", 23000),
        @("N57", -24640),
        @("H132", 2453.8948),
        @("I132", 2122.9285),
        @("K132", 6368.7855),
        @("M132", -3838.7855)
    )
    "LTW" = @(
        @("H7", 2361),
        @("I7", 2333.8333),
        @("K7", 2333.8333),
        @("M7", -2221.8333),
        @("H22", 976.6667),
        @("I22", 730),
        @("J22", 1100),
        @("K22", 730),
        @("L22", 1100),
        @("M22", -435),
        @("N22", -1690),
        @("H27", 976.6667),
        @("I27", 730),
        @("J27", 1100),
        @("K27", 730),
        @("L27", 1100),
        @("M27", -623),
        @("N27", -1314),
        @("H35", 1894),
        @("I35", 1894),
        @("K35", 1894),
        @("M35", -1558),
        @("H122", 31252488),
        @("I122", 50002260),
        @("J122", 2868.3333),
        @("K122", 150006780),
        @("L122", 8604.999899999999),
        @("M122", -150004330),
        @("N122", -13504.9999),
        @("H126", 2361),
        @("I126", 2333.8333),
        @("K126", 7001.499899999999),
        @("M126", -4531.499899999999)
    )
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $changes[$sheetName]) {
        $ref = $entry[0]
        $val = $entry[1]
        if ($null -eq $val) {
            $ws.Range($ref).ClearContents()
        } else {
            $ws.Range($ref).Value = $val
        }
    }
}

$wb.Save()
